# Bai 32 - Log4j2
# Applies the edits described in the commit diff:
#  - CustomerDataProvider!A2: "Viettel_01" -> "Viettel_03"
#  - CustomerDataProvider!A3: "Viettel_02" -> "Viettel_04"
#  - Font for both named styles (Normal / Hyperlink) changed Arial -> Aptos Narrow
#  - CustomerDataProvider sheet view selection moved from D7 to A6

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginDataProvider")
$ws2 = $wb.Worksheets.Item("CustomerDataProvider")

# --- Data edits -----------------------------------------------------------
$ws2.Range("A2").Value = "Viettel_03"
$ws2.Range("A3").Value = "Viettel_04"

# --- Font change (Arial -> Aptos Narrow) -----------------------------------
# The workbook's two fonts (Normal body text + Hyperlink) both changed from
# Arial to Aptos Narrow, so apply it across every used cell on both sheets.
$ws1.Cells.Font.Name = "Aptos Narrow"
$ws2.Cells.Font.Name = "Aptos Narrow"

# --- Selection change on CustomerDataProvider (D7 -> A6) -------------------
$ws2.Activate() | Out-Null
$ws2.Range("A6").Select() | Out-Null
